$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: plain row, same format as row 21 (B/C time values, D wrapped text)
$ws.Range("B21:D21").Copy()
$ws.Range("B22:D22").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B22").Value = 0.70833333333333337
$ws.Range("C22").Value = 0.72222222222222221
$ws.Range("D22").Value = "Character, tried to create basic sit animation inside Unity."

# Row 23: taller row (ht=30), same format as row 18
$ws.Range("B18:D18").Copy()
$ws.Range("B23:D23").PasteSpecial(-4122) # xlPasteFormats
$ws.Rows.Item(23).RowHeight = $ws.Rows.Item(18).RowHeight
$ws.Range("B23").Value = 0.72222222222222221
$ws.Range("C23").Value = 0.76388888888888884
$ws.Range("D23").Value = "Assigned materials to the car in blender, fixed the import issues."

# Row 24: plain row, same format as row 21
$ws.Range("B21:D21").Copy()
$ws.Range("B24:D24").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B24").Value = 0.76736111111111116
$ws.Range("C24").Value = 0.77569444444444446
$ws.Range("D24").Value = "Animating character in Blender."

# Row 25: plain row, same format as row 21
$ws.Range("B21:D21").Copy()
$ws.Range("B25:D25").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B25").Value = 0.77777777777777779
$ws.Range("C25").Value = 0.79583333333333339
$ws.Range("D25").Value = "Fixing up bone names in Blender…"

# Row 26: plain row, same format as row 21
$ws.Range("B21:D21").Copy()
$ws.Range("B26:D26").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B26").Value = 0.79583333333333339
$ws.Range("C26").Value = 0.80555555555555547
$ws.Range("D26").Value = "Character finally proparly imported…"

$excel.CutCopyMode = 0

# Update the view: scrolled to show the new rows, selection on D28
$excel.ActiveWindow.ScrollRow = 19
$null = $ws.Range("D28").Select()
